# "fixing test data for multi reg"
# - add a new note value for item 6 ("notes on description 6") into the
#   "testreg4" sheet at F7 (creates a new shared string entry)
# - update the remembered cell selection on both sheets (registerinfo -> B5,
#   testreg4 -> F7), leaving testreg4 as the tab-selected (active) sheet,
#   matching the order the original author last touched the sheets in.

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("registerinfo")
$wsReg4 = $wb.Worksheets.Item("testreg4")

# New note for item 6 (row 7) - this mints a new shared-string entry.
$wsReg4.Range("F7").Value = "notes on description 6"

# Restore the on-screen selections that were captured when the file was
# last saved. Select registerinfo first so that testreg4 ends up as the
# active / tab-selected sheet, same as before the edit.
$wsInfo.Range("B5").Select()
$wsReg4.Range("F7").Select()
